$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update H2 value: "Andamento" -> "Concluido"
$ws.Range("H2").Value = "Concluido"

# Add new row 3 with data
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2025-02-06"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = 34
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "2025-02-06"
$ws.Range("G3").ClearFormats()
$ws.Range("H3").Value = "Concluido"
